$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-20 (the "Application" group) had an implicit leading "Application"
# field value that was missing from column C. Shift the existing field-path
# values (columns C:F) one column to the right, and insert "Application"
# into column C for each of these rows.
for ($r = 2; $r -le 20; $r++) {
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 6).Value2 = $e
    $ws.Cells.Item($r, 5).Value2 = $d
    $ws.Cells.Item($r, 4).Value2 = $c
    $ws.Cells.Item($r, 3).Value2 = "Application"
}
